# Auto-generated edits applying the diff to Halicarnassus_Profits workbook
# Updates currentAveragePrice / LevePrice / LeveProfit computed columns (H-N)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1859.3125
$ws.Range("J17").Value = 1859.3125
$ws.Range("L17").Value = 5577.9375
$ws.Range("N17").Value = -5913.9375
$ws.Range("H86").Value = 8268.1
$ws.Range("I86").Value = 7114.8335
$ws.Range("J86").Value = 9998
$ws.Range("K86").Value = 7114.8335
$ws.Range("L86").Value = 9998
$ws.Range("M86").Value = -5991.8335
$ws.Range("N86").Value = -12244
$ws.Range("H89").Value = 8268.1
$ws.Range("I89").Value = 7114.8335
$ws.Range("J89").Value = 9998
$ws.Range("K89").Value = 35574.1675
$ws.Range("L89").Value = 49990
$ws.Range("M89").Value = -29958.1675
$ws.Range("N89").Value = -61222
$ws.Range("H96").Value = 1901
$ws.Range("I96").Value = 450
$ws.Range("J96").Value = 2191.2
$ws.Range("K96").Value = 1350
$ws.Range("L96").Value = 6573.599999999999
$ws.Range("M96").Value = 23
$ws.Range("N96").Value = -9319.599999999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1730.8
$ws.Range("I2").Value = 921.44446
$ws.Range("K2").Value = 921.44446
$ws.Range("M2").Value = -808.44446
$ws.Range("H5").Value = 155.5
$ws.Range("I5").Value = 133.5
$ws.Range("J5").Value = 177.5
$ws.Range("K5").Value = 133.5
$ws.Range("L5").Value = 177.5
$ws.Range("M5").Value = -21.5
$ws.Range("N5").Value = -401.5
$ws.Range("H39").Value = 6599.5
$ws.Range("I39").Value = 6599.5
$ws.Range("K39").Value = 6599.5
$ws.Range("M39").Value = -6079.5
$ws.Range("H45").Value = 1988.8572
$ws.Range("J45").Value = 3150.8
$ws.Range("L45").Value = 3150.8
$ws.Range("N45").Value = -3904.8
$ws.Range("H63").Value = 5649.5557
$ws.Range("I63").Value = 1549.8572
$ws.Range("J63").Value = 19998.5
$ws.Range("K63").Value = 1549.8572
$ws.Range("L63").Value = 19998.5
$ws.Range("M63").Value = -863.8571999999999
$ws.Range("N63").Value = -21370.5
$ws.Range("H66").Value = 5649.5557
$ws.Range("I66").Value = 1549.8572
$ws.Range("J66").Value = 19998.5
$ws.Range("K66").Value = 7749.286
$ws.Range("L66").Value = 99992.5
$ws.Range("M66").Value = -4317.286
$ws.Range("N66").Value = -106856.5
$ws.Range("H74").Value = 5929.6924
$ws.Range("I74").Value = 5935.273
$ws.Range("K74").Value = 5935.273
$ws.Range("M74").Value = -5061.273
$ws.Range("H77").Value = 5929.6924
$ws.Range("I77").Value = 5935.273
$ws.Range("K77").Value = 29676.365
$ws.Range("M77").Value = -25308.365
$ws.Range("H97").Value = 548.6923
$ws.Range("I97").Value = 572.75
$ws.Range("K97").Value = 572.75
$ws.Range("M97").Value = -76.75
$ws.Range("H102").Value = 2827.7058
$ws.Range("J102").Value = 10000
$ws.Range("L102").Value = 10000
$ws.Range("N102").Value = -13244
$ws.Range("H116").Value = 1730.8
$ws.Range("I116").Value = 921.44446
$ws.Range("K116").Value = 921.44446
$ws.Range("M116").Value = 1372.55554
$ws.Range("H119").Value = 60000
$ws.Range("J119").Value = 60000
$ws.Range("L119").Value = 60000
$ws.Range("N119").Value = -69676
$ws.Range("H132").Value = 1432.6666
$ws.Range("I132").Value = 1432.6666
$ws.Range("K132").Value = 4297.9998
$ws.Range("M132").Value = -1767.9998

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1730.8
$ws.Range("I3").Value = 921.44446
$ws.Range("K3").Value = 921.44446
$ws.Range("M3").Value = -807.44446
$ws.Range("H4").Value = 155.5
$ws.Range("I4").Value = 133.5
$ws.Range("J4").Value = 177.5
$ws.Range("K4").Value = 133.5
$ws.Range("L4").Value = 177.5
$ws.Range("M4").Value = -18.5
$ws.Range("N4").Value = -407.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H55").Value = 40000
$ws.Range("I55").Value = 40000
$ws.Range("K55").Value = 40000
$ws.Range("M55").Value = -39685
$ws.Range("H58").Value = 4642.643
$ws.Range("I58").Value = 4205.2
$ws.Range("K58").Value = 4205.2
$ws.Range("M58").Value = -4002.2
$ws.Range("H88").Value = 14665.833
$ws.Range("J88").Value = 14665.833
$ws.Range("L88").Value = 14665.833
$ws.Range("N88").Value = -15477.833
$ws.Range("H91").Value = 14665.833
$ws.Range("J91").Value = 14665.833
$ws.Range("L91").Value = 14665.833
$ws.Range("N91").Value = -17473.833
$ws.Range("H132").Value = 1593.1538
$ws.Range("I132").Value = 1673.8182
$ws.Range("K132").Value = 5021.4546
$ws.Range("M132").Value = -2491.4546
$ws.Range("H136").Value = 4642.643
$ws.Range("I136").Value = 4205.2
$ws.Range("K136").Value = 12615.6
$ws.Range("M136").Value = -10065.6
$ws.Range("H140").Value = 57890
$ws.Range("J140").Value = 57890
$ws.Range("L140").Value = 57890
$ws.Range("N140").Value = -68250

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 294
$ws.Range("J22").Value = 295.55554
$ws.Range("L22").Value = 886.66662
$ws.Range("N22").Value = -1224.66662
$ws.Range("H27").Value = 294
$ws.Range("J27").Value = 295.55554
$ws.Range("L27").Value = 886.66662
$ws.Range("N27").Value = -1090.66662
$ws.Range("H57").Value = 1565.8889
$ws.Range("J57").Value = 2500
$ws.Range("L57").Value = 7500
$ws.Range("N57").Value = -8618
$ws.Range("H60").Value = 2111.818
$ws.Range("I60").Value = 250
$ws.Range("J60").Value = 2298
$ws.Range("K60").Value = 750
$ws.Range("L60").Value = 6894
$ws.Range("M60").Value = -499
$ws.Range("N60").Value = -7396
$ws.Range("H131").Value = 1720.7
$ws.Range("I131").Value = 653
$ws.Range("J131").Value = 2178.2856
$ws.Range("K131").Value = 1959
$ws.Range("L131").Value = 6534.8568
$ws.Range("M131").Value = 3081
$ws.Range("N131").Value = -16614.8568
$ws.Range("H136").Value = 5500
$ws.Range("I136").Value = 2250
$ws.Range("K136").Value = 6750
$ws.Range("M136").Value = -1650

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1964.4546
$ws.Range("I80").Value = 1800.4286
$ws.Range("J80").Value = 2251.5
$ws.Range("K80").Value = 1800.4286
$ws.Range("L80").Value = 2251.5
$ws.Range("M80").Value = -802.4286
$ws.Range("N80").Value = -4247.5
$ws.Range("H83").Value = 1964.4546
$ws.Range("I83").Value = 1800.4286
$ws.Range("J83").Value = 2251.5
$ws.Range("K83").Value = 9002.143
$ws.Range("L83").Value = 11257.5
$ws.Range("M83").Value = -4010.143
$ws.Range("N83").Value = -21241.5
$ws.Range("H132").Value = 3126.75
$ws.Range("I132").Value = 2859.1428
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 8577.4284
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -6047.428400000001
$ws.Range("N132").Value = -20060

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2001
$ws.Range("I16").Value = 2001
$ws.Range("K16").Value = 2001
$ws.Range("M16").Value = -1831
$ws.Range("H61").Value = 5011
$ws.Range("I61").Value = 2026.3334
$ws.Range("K61").Value = 2026.3334
$ws.Range("M61").Value = -1824.3334
$ws.Range("H113").Value = 5011
$ws.Range("I113").Value = 2026.3334
$ws.Range("K113").Value = 2026.3334
$ws.Range("M113").Value = 143.6666

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7675.6924
$ws.Range("J62").Value = 8333.143
$ws.Range("L62").Value = 8333.143
$ws.Range("N62").Value = -9581.143
$ws.Range("H65").Value = 7675.6924
$ws.Range("J65").Value = 8333.143
$ws.Range("L65").Value = 41665.715
$ws.Range("N65").Value = -47905.715
$ws.Range("H136").Value = 3470.9167
$ws.Range("J136").Value = 5281.3
$ws.Range("L136").Value = 15843.9
$ws.Range("N136").Value = -20943.9
